$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force the cell to be written as text even when the string looks
    # like a number (e.g. "20020808"), then strip the temporary
    # number-format override again so the cell keeps the default style.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# ---- Row 2: replace existing patent record with new data ----
$ws.Range("B2").Value = "US20020107833A1"
$ws.Range("C2").Value = "Method and system for tracking equipment usage information"
Set-TextValue $ws "D2" "22585439"
Set-TextValue $ws "E2" "20011116"
$ws.Range("F2").Value = 19991029
Set-TextValue $ws "G2" "20020808"
$ws.Range("H2").Value = "KERKINNI FUAT J"
$ws.Range("I2").Value = "KERKINNI FUAT"
$ws.Range("J2").Value = "US2002107833A1| AU2619801A| WO0135679A3| WO0135679A2"

# ---- Row 3: new patent record ----
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "JP2002196817A"
$ws.Range("C3").Value = "SYSTEM FOR MANAGING WORK VEHICLE"
Set-TextValue $ws "D3" "18859294"
Set-TextValue $ws "E3" "20001225"
$ws.Range("F3").Value = 20001225
Set-TextValue $ws "G3" "20020712"
$ws.Range("H3").Value = "ISEKI & CO LTD"
$ws.Range("I3").Value = "IKEUCHI NOBUAKI, WATABE TOMOAKI, SHINODA MASANORI, SAKATA CHIKANO"
$ws.Range("J3").Value = "JP2002196817A"

# ---- Row 4: new patent record ----
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "US6434512B1"
$ws.Range("C4").Value = "Modular data collection and analysis system"
Set-TextValue $ws "D4" "27489613"
Set-TextValue $ws "E4" "19990930"
$ws.Range("F4").Value = 19980402
Set-TextValue $ws "G4" "20020813"
$ws.Range("H4").Value = "RELIANCE ELECTRIC TECHNOLOGIES, LLC"
$ws.Range("I4").Value = "DISCENZO FREDERICK"
$ws.Range("J4").Value = "US7690246B1| US6546785B1| US7493799B1| US6877360B1| US6286363B1| US6434512B1| US6023961A| US6196057B1| US7134323B1| US6295510B1| US6324899B1"

# ---- Row 5: new patent record ----
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = " US20020111725A1"
$ws.Range("C5").Value = "Method and apparatus for risk-related use of vehicle communication system data"
Set-TextValue $ws "D5" "27575201"
Set-TextValue $ws "E5" "20010716"
$ws.Range("F5").Value = 20000717
Set-TextValue $ws "G5" "20020815"
$ws.Range("H5").Value = "BURGE JOHN R"
$ws.Range("I5").Value = "BURGE JOHN"
$ws.Range("J5").Value = "US2002103622A1| US2002111725A1"

# ---- Copy the formatting from A2 (bordered, bold, centered) down to the new A column cells ----
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
